$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D width: raw OOXML width=30 (ColumnWidth value tuned so the stored width comes out to exactly 30)
$ws.Columns.Item(4).ColumnWidth = 29.17

# Drop every existing hyperlink on the sheet so the relationships can be rebuilt fresh, in the final row order
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Range('A2').Value = '2025-11-27 12:37:57'
$ws.Range('B2').Value = '【急募】pythonのコードのMac環境用インストーラー作成(Windows版は作成済)'
$ws.Range('C2').Value = 'システム開発'
$ws.Range('D2').Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range('E2').Value = '期限情報なし'
$ws.Range('F2').Value = 'https://www.lancers.jp/work/detail/5442448'
$ws.Range('G2').Value = 190
$ws.Range('H2').Value = '🔥Python'

# Row 3
$ws.Range('A3').Value = '2025-11-27 12:37:57'
$ws.Range('B3').Value = 'Threads・X対応|スプレッドシート連携の自動投稿ツール開発をお願いできる方'
$ws.Range('C3').Value = 'システム開発'
$ws.Range('D3').Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range('E3').Value = '期限情報なし'
$ws.Range('F3').Value = 'https://www.lancers.jp/work/detail/5442360'
$ws.Range('G3').Value = 120
$ws.Range('H3').Value = '◆ツール,開発'

# Row 4
$ws.Range('A4').Value = '2025-11-27 12:37:57'
$ws.Range('B4').Value = '【Java/Tomcat】スクラッチ構築の予約サイトにおける複数バグの修正依頼'
$ws.Range('C4').Value = 'システム開発'
$ws.Range('D4').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E4').Value = '期限情報なし'
$ws.Range('F4').Value = 'https://www.lancers.jp/work/detail/5442482'
$ws.Range('G4').Value = 103
$ws.Range('H4').Value = '★Java ◇サイト'

# Row 5
$ws.Range('A5').Value = '2025-11-27 12:37:57'
$ws.Range('B5').Value = '【心理学実験】Javaシステム改修と機能拡張の依頼'
$ws.Range('C5').Value = 'システム開発'
$ws.Range('D5').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E5').Value = '期限情報なし'
$ws.Range('F5').Value = 'https://www.lancers.jp/work/detail/5442416'
$ws.Range('G5').Value = 100
$ws.Range('H5').Value = '★Java'

# Row 6
$ws.Range('A6').Value = '2025-11-27 12:37:57'
$ws.Range('B6').Value = 'Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Range('C6').Value = 'システム開発'
$ws.Range('D6').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E6').Value = '期限情報なし'
$ws.Range('F6').Value = 'https://www.lancers.jp/work/detail/5442063'
$ws.Range('G6').Value = 85
$ws.Range('H6').Value = '★Java'

# Row 7
$ws.Range('A7').Value = '2025-11-27 12:37:57'
$ws.Range('B7').Value = 'クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Range('C7').Value = 'システム開発'
$ws.Range('D7').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E7').Value = '期限情報なし'
$ws.Range('F7').Value = 'https://www.lancers.jp/work/detail/5442064'
$ws.Range('G7').Value = 38
$ws.Range('H7').Value = '◇管理'

# Row 8
$ws.Range('A8').Value = '2025-11-27 12:37:57'
$ws.Range('B8').Value = '【急募】運用中のshopifyサイトにGMOイプシロン決済導入'
$ws.Range('C8').Value = 'システム開発'
$ws.Range('D8').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E8').Value = '期限情報なし'
$ws.Range('F8').Value = 'https://www.lancers.jp/work/detail/5442625'
$ws.Range('G8').Value = 33
$ws.Range('H8').Value = '◇サイト'

# Row 9
$ws.Range('A9').Value = '2025-11-27 12:37:57'
$ws.Range('B9').Value = 'Access業務システムのクラウド化(ZOHO Creator使用)をお手伝いください!(再依頼)'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('F9').Value = 'https://www.lancers.jp/work/detail/5442153'
$ws.Range('G9').Value = 40

# Row 10
$ws.Range('A10').Value = '2025-11-27 12:37:57'
$ws.Range('B10').Value = '【募集】Amazonフラットファイル(ブラウズノード検証)'
$ws.Range('C10').Value = 'システム開発'
$ws.Range('D10').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E10').Value = '期限情報なし'
$ws.Range('F10').Value = 'https://www.lancers.jp/work/detail/5442106'
$ws.Range('G10').Value = 13

# Row 11
$ws.Range('A11').Value = '2025-11-27 12:37:57'
$ws.Range('B11').Value = '【急募】ex4ファイルをmq4ファイルに変換していただける方'
$ws.Range('C11').Value = 'システム開発'
$ws.Range('D11').Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range('E11').Value = '期限情報なし'
$ws.Range('F11').Value = 'https://www.lancers.jp/work/detail/5442432'
$ws.Range('G11').Value = 10

# Row 12
$ws.Range('A12').Value = '2025-11-27 12:37:57'
$ws.Range('B12').Value = '【急募】ex4ファイルをmq4ファイルに変換していただける方'
$ws.Range('C12').Value = 'システム開発'
$ws.Range('D12').Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range('E12').Value = '期限情報なし'
$ws.Range('F12').Value = 'https://www.lancers.jp/work/detail/5442169'
$ws.Range('G12').Value = 10

# Rebuild hyperlinks for column F in row order so relationship ids come out sequential (rId1..rId11)
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5442448')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5442360')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5442482')
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5442416')
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5442063')
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5442064')
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5442625')
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5442153')
$ws.Hyperlinks.Add($ws.Range('F10'), 'https://www.lancers.jp/work/detail/5442106')
$ws.Hyperlinks.Add($ws.Range('F11'), 'https://www.lancers.jp/work/detail/5442432')
$ws.Hyperlinks.Add($ws.Range('F12'), 'https://www.lancers.jp/work/detail/5442169')

# Re-apply the built-in Hyperlink cell style so all link cells share the same style index as before
$ws.Range('F2:F12').Style = 'Hyperlink'
